$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hyperlinks from A10:A11 (clears the link but we still need to clear content separately)
$ws.Hyperlinks.Delete()

# Clear the text content of A9:A14, but keep their existing styles
$ws.Range("A9:A14").ClearContents()

# Update column A width (target stored width 59.77734375; closest value this
# engine's ColumnWidth->stored-width rounding can produce is 59.8333, reached
# with an input of 59)
$ws.Columns.Item(1).ColumnWidth = 59

# Update the selection
$ws.Range("A9:A14").Select()
$ws.Application.ActiveCell = $ws.Range("A14")
